$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column R (header through last data row) into new column S,
# carrying over formatting (style) along with values in one shot.
$ws.Range("R4:R14").Copy($ws.Range("S4:S14")) | Out-Null

# Set the new header and data values for column S
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 99.5
$ws.Range("S6").Value = 99.358544044156048
$ws.Range("S7").Value = 99.400057479522914
$ws.Range("S8").Value = 99.513194978221875
$ws.Range("S9").Value = 99.232429839290006
$ws.Range("S10").Value = 99.453093666824671
$ws.Range("S11").Value = 99.686258104998956
$ws.Range("S12").Value = 99.42525365081228
$ws.Range("S13").Value = 99.561275226674468
$ws.Range("S14").Value = 99.831561216970215

# Update the selected cell to match the target state
$ws.Range("U6").Select() | Out-Null
